$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(4, "China", 81171, 78, 73159, 4735, 1573, 7, 3277)
    ,@(5, "Italia", 63927, 0, 7432, 50418, 3204, 0, 6077)
    ,@(6, "Estados Unidos", 46168, 2434, 295, 45291, 1040, 29, 582)
    ,@(7, "España", 39673, 4537, 3794, 33183, 2355, 385, 2696)
    ,@(8, "Alemania", 30150, 1094, 453, 29567, 23, 7, 130)
    ,@(9, "Iran", 24811, 1762, 8913, 13964, 0, 122, 1934)
    ,@(10, "Francia", 19856, 0, 2200, 16796, 2082, 0, 860)
    ,@(11, "Suiza", 9117, 322, 131, 8864, 141, 2, 122)
    ,@(12, "Corea del Sur", 9037, 76, 3507, 5410, 59, 9, 120)
    ,@(13, "Reino Unido", 6650, 0, 135, 6180, 20, 0, 335)
    ,@(14, "Austria", 4791, 317, 9, 4757, 17, 4, 25)
    ,@(15, "Paises Bajos", 4749, 0, 2, 4534, 435, 0, 213)
    ,@(16, "Belgica", 4269, 526, 461, 3686, 381, 34, 122)
    ,@(17, "Noruega", 2647, 22, 6, 2630, 42, 1, 11)
    ,@(18, "Australia", 2136, 249, 118, 2010, 11, 1, 8)
    ,@(19, "Canada", 2091, 0, 112, 1955, 1, 0, 24)
    ,@(20, "Portugal", 2060, 0, 14, 2023, 47, 0, 23)
    ,@(21, "Suecia", 2059, 13, 16, 2010, 110, 6, 33)
    ,@(22, "Brasil", 1960, 36, 2, 1924, 18, 0, 34)
    ,@(23, "Israel", 1656, 214, 49, 1606, 31, 0, 1)
    ,@(24, "Malasia", 1624, 106, 183, 1426, 64, 1, 15)
    ,@(25, "Dinamarca", 1577, 117, 1, 1544, 69, 8, 32)
    ,@(26, "Turquia", 1529, 0, 0, 1492, 0, 0, 37)
    ,@(27, "Chequia", 1289, 53, 7, 1280, 19, 1, 2)
    ,@(28, "Japon", 1140, 12, 285, 813, 54, 0, 42)
    ,@(29, "Irlanda", 1125, 0, 5, 1114, 29, 0, 6)
    ,@(30, "Ecuador", 981, 0, 3, 960, 2, 0, 18)
    ,@(31, "Pakistan", 908, 33, 13, 888, 0, 1, 7)
    ,@(32, "Luxemburgo", 875, 0, 6, 861, 3, 0, 8)
    ,@(33, "Tailandia", 827, 106, 52, 771, 7, 3, 4)
    ,@(34, "Finlandia", 792, 92, 10, 781, 11, 0, 1)
    ,@(35, "Polonia", 774, 25, 1, 764, 3, 1, 9)
    ,@(36, "Rumania", 762, 186, 79, 675, 15, 1, 8)
    ,@(37, "Chile", 746, 0, 11, 733, 7, 0, 2)
    ,@(38, "Crucero", 712, 0, 587, 115, 15, 2, 10)
    ,@(39, "Grecia", 695, 0, 29, 649, 35, 0, 17)
    ,@(40, "Indonesia", 686, 107, 30, 601, 0, 6, 55)
    ,@(41, "Islandia", 588, 0, 51, 535, 14, 1, 2)
    ,@(42, "Arabia Saudita", 562, 0, 19, 543, 0, 0, 0)
    ,@(43, "Sudafrica", 554, 152, 4, 550, 2, 0, 0)
    ,@(44, "Filipinas", 552, 90, 20, 497, 1, 2, 35)
    ,@(45, "India", 511, 12, 37, 464, 0, 0, 10)
    ,@(46, "Singapur", 509, 0, 152, 355, 14, 0, 2)
    ,@(47, "Catar", 501, 0, 37, 464, 6, 0, 0)
    ,@(48, "Rusia", 495, 57, 22, 472, 0, 0, 1)
    ,@(49, "Eslovenia", 442, 0, 0, 439, 12, 0, 3)
    ,@(50, "Peru", 395, 0, 1, 389, 19, 0, 5)
    ,@(51, "Barein", 390, 13, 164, 224, 3, 0, 2)
    ,@(52, "Hong Kong", 386, 29, 102, 280, 4, 0, 4)
    ,@(53, "Estonia", 369, 17, 7, 362, 4, 0, 0)
    ,@(54, "Mexico", 367, 51, 4, 359, 1, 1, 4)
    ,@(55, "Egipto", 366, 0, 68, 279, 0, 0, 19)
    ,@(56, "Croacia", 361, 46, 5, 355, 6, 0, 1)
    ,@(57, "Panama", 345, 0, 1, 338, 33, 0, 6)
    ,@(58, "Colombia", 306, 29, 6, 297, 0, 0, 3)
    ,@(59, "Libano", 304, 37, 8, 292, 4, 0, 4)
    ,@(60, "Argentina", 301, 0, 51, 246, 0, 0, 4)
    ,@(61, "Irak", 266, 0, 62, 181, 0, 0, 23)
    ,@(62, "Serbia", 249, 0, 3, 243, 12, 0, 3)
    ,@(63, "Republica Dominicana", 245, 0, 3, 239, 0, 0, 3)
    ,@(64, "Armenia", 235, 0, 2, 233, 6, 0, 0)
    ,@(65, "Argelia", 230, 0, 65, 148, 0, 0, 17)
    ,@(66, "Taiwan", 216, 21, 29, 185, 0, 0, 2)
    ,@(67, "Eslovaquia", 204, 18, 7, 197, 2, 0, 0)
    ,@(68, "Bulgaria", 202, 1, 3, 196, 8, 0, 3)
    ,@(69, "Emiratos Arabes Unidos", 198, 0, 41, 155, 2, 0, 2)
    ,@(70, "Letonia", 197, 17, 1, 196, 0, 0, 0)
    ,@(71, "Kuwait", 191, 2, 39, 152, 5, 0, 0)
    ,@(72, "Lituania", 187, 8, 1, 185, 1, 0, 1)
    ,@(73, "San Marino", 187, 0, 4, 163, 13, 0, 20)
    ,@(74, "Hungria", 187, 20, 21, 157, 6, 1, 9)
    ,@(75, "Principado de Andorra", 164, 31, 1, 162, 7, 0, 1)
    ,@(76, "Uruguay", 162, 0, 0, 162, 3, 0, 0)
    ,@(77, "Costa Rica", 158, 0, 2, 154, 2, 0, 2)
    ,@(78, "Nueva Zelanda", 155, 53, 12, 143, 0, 0, 0)
    ,@(79, "Bosnia y Herzegovina", 150, 14, 2, 147, 1, 0, 1)
    ,@(80, "Marruecos", 143, 0, 5, 134, 1, 0, 4)
    ,@(81, "Republica de Macedonia", 136, 0, 1, 133, 1, 0, 2)
    ,@(82, "Jordania", 127, 0, 1, 126, 0, 0, 0)
    ,@(83, "Albania", 123, 19, 10, 108, 2, 1, 5)
    ,@(84, "Vietnam", 123, 0, 17, 106, 3, 0, 0)
    ,@(85, "Islas Feroe", 122, 4, 23, 99, 0, 0, 0)
    ,@(86, "Republica de Chipre", 116, 0, 3, 112, 3, 0, 1)
    ,@(87, "Burkina Faso", 114, 15, 7, 103, 0, 0, 4)
    ,@(88, "Malta", 110, 3, 2, 108, 1, 0, 0)
    ,@(89, "Moldavia", 109, 0, 2, 106, 10, 0, 1)
    ,@(90, "Brunei", 104, 13, 2, 102, 2, 0, 0)
    ,@(91, "Sri Lanka", 100, 3, 2, 98, 2, 0, 0)
    ,@(92, "Ucrania", 97, 24, 1, 93, 0, 0, 3)
    ,@(93, "Tunez", 90, 1, 1, 86, 11, 0, 3)
    ,@(94, "Camboya", 87, 0, 2, 85, 1, 0, 0)
    ,@(95, "Senegal", 86, 7, 8, 78, 0, 0, 0)
    ,@(96, "Venezuela", 84, 0, 15, 69, 2, 0, 0)
    ,@(97, "Oman", 84, 18, 17, 67, 0, 0, 0)
    ,@(98, "Bielorrusia", 81, 0, 22, 59, 0, 0, 0)
    ,@(99, "Reunion", 75, 4, 1, 74, 0, 0, 0)
    ,@(100, "Azerbaiyan", 72, 0, 10, 61, 3, 0, 1)
    ,@(101, "Kazajistan", 68, 6, 0, 68, 0, 0, 0)
    ,@(102, "Georgia", 67, 6, 9, 58, 1, 0, 0)
    ,@(103, "Camerun", 66, 10, 2, 64, 0, 0, 0)
    ,@(104, "Guadalupe", 62, 0, 0, 61, 4, 0, 1)
    ,@(105, "Estado de Palestina", 60, 1, 17, 43, 0, 0, 0)
    ,@(106, "Martinica", 53, 0, 0, 52, 7, 0, 1)
    ,@(107, "Trinidad yTobago", 52, 1, 0, 52, 0, 0, 0)
    ,@(108, "Ghana", 52, 25, 0, 50, 0, 0, 2)
    ,@(109, "Liechtenstein", 51, 0, 0, 51, 0, 0, 0)
    ,@(110, "Uzbekistan", 50, 4, 0, 50, 0, 0, 0)
    ,@(111, "Consejo Danes para los Refugiados", 45, 9, 0, 43, 0, 0, 2)
    ,@(112, "Kirguistan", 42, 26, 0, 42, 0, 0, 0)
    ,@(113, "Afganistan", 42, 2, 1, 40, 0, 0, 1)
    ,@(114, "Cuba", 40, 0, 0, 39, 3, 0, 1)
    ,@(115, "Nigeria", 40, 0, 2, 37, 0, 0, 1)
    ,@(116, "Puerto Rico", 39, 8, 1, 36, 0, 0, 2)
    ,@(117, "Banglades", 39, 6, 5, 30, 0, 1, 4)
    ,@(118, "Mayotte", 36, 12, 0, 36, 0, 0, 0)
    ,@(119, "Ruanda", 36, 0, 0, 36, 0, 0, 0)
    ,@(120, "Mauricio", 36, 0, 0, 34, 1, 0, 2)
    ,@(121, "Honduras", 30, 0, 0, 30, 0, 0, 0)
    ,@(122, "Guam", 29, 0, 0, 28, 0, 0, 1)
    ,@(123, "Montenegro", 29, 2, 0, 28, 0, 0, 1)
    ,@(124, "Bolivia", 28, 1, 0, 28, 0, 0, 0)
    ,@(125, "Paraguay", 27, 5, 0, 25, 1, 1, 2)
    ,@(126, "Costa de Marfil", 25, 0, 2, 23, 0, 0, 0)
    ,@(127, "Macao", 25, 0, 10, 15, 0, 0, 0)
    ,@(128, "Polinesia Francesa", 23, 5, 0, 23, 0, 0, 0)
    ,@(129, "Monaco", 23, 0, 1, 22, 0, 0, 0)
    ,@(130, "Guayana Francesa", 23, 3, 6, 17, 0, 0, 0)
    ,@(131, "Jamaica", 21, 2, 2, 18, 0, 0, 1)
    ,@(132, "Isla de Man", 20, 7, 0, 20, 0, 0, 0)
    ,@(133, "Guatemala", 20, 0, 0, 19, 0, 0, 1)
    ,@(134, "Guyana", 20, 0, 0, 19, 0, 0, 1)
    ,@(135, "Togo", 18, 0, 0, 18, 0, 0, 0)
    ,@(136, "Barbados", 17, 0, 0, 17, 0, 0, 0)
    ,@(137, "Madagascar", 17, 5, 0, 17, 0, 0, 0)
    ,@(138, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
    ,@(139, "Kenia", 16, 0, 0, 16, 0, 0, 0)
    ,@(140, "Gibraltar", 15, 0, 5, 10, 0, 0, 0)
    ,@(141, "Maldivas", 13, 0, 5, 8, 0, 0, 0)
    ,@(142, "Etiopia", 12, 1, 0, 12, 0, 0, 0)
    ,@(143, "Tanzania", 12, 0, 0, 12, 0, 0, 0)
    ,@(144, "Aruba", 12, 3, 1, 11, 0, 0, 0)
    ,@(145, "Nueva Caledonia", 10, 2, 0, 10, 0, 0, 0)
    ,@(146, "Mongolia", 10, 0, 0, 10, 0, 0, 0)
    ,@(147, "Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0)
    ,@(148, "Uganda", 9, 0, 0, 9, 0, 0, 0)
    ,@(149, "San Martin (Parte Francesa)", 8, 0, 0, 8, 0, 0, 0)
    ,@(150, "Seychelles", 7, 0, 0, 7, 0, 0, 0)
    ,@(151, "Bermudas", 6, 0, 0, 6, 0, 0, 0)
    ,@(152, "Surinam", 6, 1, 0, 6, 0, 0, 0)
    ,@(153, "Haiti", 6, 0, 0, 6, 0, 0, 0)
    ,@(154, "Benin", 6, 0, 0, 6, 0, 0, 0)
    ,@(155, "Gabon", 6, 0, 0, 5, 0, 0, 1)
    ,@(156, "El Salvador", 5, 2, 0, 5, 0, 0, 0)
    ,@(157, "Islas Caimanes", 5, 0, 0, 4, 0, 0, 1)
    ,@(158, "Bahamas", 4, 0, 0, 4, 0, 0, 0)
    ,@(159, "Suazilandia", 4, 0, 0, 4, 0, 0, 0)
    ,@(160, "Guinea", 4, 0, 0, 4, 0, 0, 0)
    ,@(161, "Fiyi", 4, 1, 0, 4, 0, 0, 0)
    ,@(162, "Groenlandia", 4, 0, 0, 4, 0, 0, 0)
    ,@(163, "Congo", 4, 0, 0, 4, 0, 0, 0)
    ,@(164, "Namibia", 4, 0, 0, 4, 0, 0, 0)
    ,@(165, "Curazao", 4, 0, 0, 3, 0, 0, 1)
    ,@(166, "Republica del Chad", 3, 1, 0, 3, 0, 0, 0)
    ,@(167, "Antigua y Barbuda", 3, 0, 0, 3, 0, 0, 0)
    ,@(168, "Liberia", 3, 0, 0, 3, 0, 0, 0)
    ,@(169, "Angola", 3, 0, 0, 3, 0, 0, 0)
    ,@(170, "Niger", 3, 0, 0, 3, 0, 0, 0)
    ,@(171, "Santa Lucia", 3, 0, 0, 3, 0, 0, 0)
    ,@(172, "Zambia", 3, 0, 0, 3, 0, 0, 0)
    ,@(173, "Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0)
    ,@(174, "Republica de Yibuti", 3, 0, 0, 3, 0, 0, 0)
    ,@(175, "San Bartolome", 3, 0, 0, 3, 0, 0, 0)
    ,@(176, "Cabo Verde", 3, 0, 0, 2, 0, 1, 1)
    ,@(177, "Zimbabue", 3, 0, 0, 2, 0, 0, 1)
    ,@(178, "Nicaragua", 2, 0, 0, 2, 0, 0, 0)
    ,@(179, "Laos", 2, 2, 0, 2, 0, 0, 0)
    ,@(180, "Butan", 2, 0, 0, 2, 0, 0, 0)
    ,@(181, "Dominica", 2, 0, 0, 2, 0, 0, 0)
    ,@(182, "Mauritania", 2, 0, 0, 2, 0, 0, 0)
    ,@(183, "San Martin (Parte Holandesa)", 2, 0, 0, 2, 0, 0, 0)
    ,@(184, "Birmania", 2, 0, 0, 2, 0, 0, 0)
    ,@(185, "Sudan", 2, 0, 0, 1, 0, 0, 1)
    ,@(186, "Nepal", 2, 0, 1, 1, 0, 0, 0)
    ,@(187, "Gambia", 2, 0, 0, 1, 0, 0, 1)
    ,@(188, "Montserrat", 1, 0, 0, 1, 0, 0, 0)
    ,@(189, "San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0)
    ,@(190, "Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
    ,@(191, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
    ,@(192, "Eritrea", 1, 0, 0, 1, 0, 0, 0)
    ,@(193, "Mozambique", 1, 0, 0, 1, 0, 0, 0)
    ,@(194, "Siria", 1, 0, 0, 1, 0, 0, 0)
    ,@(195, "Granada", 1, 0, 0, 1, 0, 0, 0)
    ,@(196, "Santa Sede", 1, 0, 0, 1, 0, 0, 0)
    ,@(197, "Somalia", 1, 0, 0, 1, 0, 0, 0)
    ,@(198, "Islas Turcas y Caicos", 1, 0, 0, 1, 0, 0, 0)
    ,@(199, "Belice", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 12:46"
